$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1035:1036, shifting existing rows (1035-1115) down to (1037-1117)
$ws.Rows("1035:1036").Insert()

# Row 1035 - new data
$ws.Cells.Item(1035, 1).Value = 10
$ws.Cells.Item(1035, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1035, 3).Value = "La Araucanía"
$ws.Cells.Item(1035, 4).Value = 45106
$ws.Cells.Item(1035, 5).Value = 9
$ws.Cells.Item(1035, 6).Value = 100112003
$ws.Cells.Item(1035, 7).Value = "Ajo"
$ws.Cells.Item(1035, 8).Value = "Chino"
$ws.Cells.Item(1035, 9).Value = "Primera"
$ws.Cells.Item(1035, 10).Value = 185
$ws.Cells.Item(1035, 11).Value = 20000
$ws.Cells.Item(1035, 12).Value = 20000
$ws.Cells.Item(1035, 13).Value = 20000
$ws.Cells.Item(1035, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(1035, 15).Value = "China"
$ws.Cells.Item(1035, 16).Value = 2000
$ws.Cells.Item(1035, 17).Value = 10
$ws.Cells.Item(1035, 18).Value = "Hortaliza"

# Row 1036 - new data
$ws.Cells.Item(1036, 1).Value = 10
$ws.Cells.Item(1036, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1036, 3).Value = "La Araucanía"
$ws.Cells.Item(1036, 4).Value = 45106
$ws.Cells.Item(1036, 5).Value = 9
$ws.Cells.Item(1036, 6).Value = 100112003
$ws.Cells.Item(1036, 7).Value = "Ajo"
$ws.Cells.Item(1036, 8).Value = "Chino"
$ws.Cells.Item(1036, 9).Value = "Primera"
$ws.Cells.Item(1036, 10).Value = 280
$ws.Cells.Item(1036, 11).Value = 22000
$ws.Cells.Item(1036, 12).Value = 22000
$ws.Cells.Item(1036, 13).Value = 22000
$ws.Cells.Item(1036, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(1036, 15).Value = "China"
$ws.Cells.Item(1036, 16).Value = 2200
$ws.Cells.Item(1036, 17).Value = 10
$ws.Cells.Item(1036, 18).Value = "Hortaliza"
